$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to be treated as text so that
# numeric-looking strings (e.g. "1.028") are not auto-converted to numbers,
# preserving their exact original text formatting.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.596.67'
$ws.Range("E2").Value = '  +2.47%  '
$ws.Range("D3").Value = '1.851.68'
$ws.Range("E3").Value = '  +2.14%  '
$ws.Range("E4").Value = '  +2.79%  '
$ws.Range("D5").Value = '321.28'
$ws.Range("E5").Value = '  +3.16%  '
$ws.Range("D6").Value = '1.028'
$ws.Range("E6").Value = '  +2.45%  '
$ws.Range("D7").Value = '0.4392'
$ws.Range("E7").Value = '  +2.49%  '
$ws.Range("D8").Value = '0.3784'
$ws.Range("E8").Value = '  +2.72%  '
$ws.Range("D9").Value = '0.07403'
$ws.Range("E9").Value = '  +2.26%  '
$ws.Range("D10").Value = '0.8756'
$ws.Range("E10").Value = '  +1.59%  '
$ws.Range("D11").Value = '21.52'
$ws.Range("E11").Value = '  +1.77%  '
$ws.Range("D12").Value = '1.859.28'
$ws.Range("E12").Value = '  -8.84%  '
$ws.Range("D13").Value = '5.517'
$ws.Range("E13").Value = '  +2.54%  '
$ws.Range("D14").Value = '6.693'
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("D15").Value = '0.07211'
$ws.Range("E15").Value = '  +4.62%  '
$ws.Range("D16").Value = '83.10'
$ws.Range("E16").Value = '  +3.10%  '
$ws.Range("D17").Value = '1.034'
$ws.Range("E17").Value = '  +2.52%  '
$ws.Range("D18").Value = '0.000009050'
$ws.Range("E18").Value = '  +2.21%  '
$ws.Range("D19").Value = '1.029'
$ws.Range("E19").Value = '  +2.49%  '
$ws.Range("D20").Value = '15.41'
$ws.Range("E20").Value = '  +1.56%  '
$ws.Range("D21").Value = '27.607.15'
$ws.Range("E21").Value = '  +2.37%  '
$ws.Range("D22").Value = '5.256'
$ws.Range("D23").Value = '11.37'
$ws.Range("E23").Value = '  +3.22%  '
$ws.Range("D24").Value = '157.88'
$ws.Range("E24").Value = '  +2.67%  '
$ws.Range("D25").Value = '1.919'
$ws.Range("E25").Value = '  +1.82%  '
$ws.Range("D26").Value = '18.74'
$ws.Range("E26").Value = '  +2.65%  '
$ws.Range("D27").Value = '1.970'
$ws.Range("E27").Value = '  +4.23%  '
$ws.Range("D28").Value = '5.276'
$ws.Range("E28").Value = '  +1.34%  '
$ws.Range("D29").Value = '117.07'
$ws.Range("E29").Value = '  +1.93%  '
$ws.Range("D30").Value = '0.09055'
$ws.Range("E30").Value = '  +1.35%  '
$ws.Range("D31").Value = '1.198'
$ws.Range("E31").Value = '  +3.10%  '
$ws.Range("D32").Value = '0.7601'
$ws.Range("E32").Value = '  +2.60%  '
$ws.Range("D33").Value = '4.520'
$ws.Range("E33").Value = '  +2.07%  '
$ws.Range("D34").Value = '2.882'
$ws.Range("E34").Value = '  +3.14%  '
$ws.Range("D35").Value = '1.029'
$ws.Range("E35").Value = '  +2.06%  '
$ws.Range("D36").Value = '1.150'
$ws.Range("E36").Value = '  +3.11%  '
$ws.Range("D37").Value = '0.01973'
$ws.Range("E37").Value = '  +2.50%  '
$ws.Range("D38").Value = '0.05309'
$ws.Range("E38").Value = '  +1.88%  '
$ws.Range("D39").Value = '2.819'
$ws.Range("E39").Value = '  +2.41%  '
$ws.Range("D40").Value = '0.5151'
$ws.Range("E40").Value = '  +1.37%  '
$ws.Range("D41").Value = '0.1675'
$ws.Range("E41").Value = '  +1.96%  '
$ws.Range("D42").Value = '6.755'
$ws.Range("E42").Value = '  +4.94%  '
$ws.Range("D43").Value = '8.477'
$ws.Range("E43").Value = '  +2.69%  '
$ws.Range("D44").Value = '108.69'
$ws.Range("E44").Value = '  +1.79%  '
$ws.Range("D45").Value = '10.57'
$ws.Range("E45").Value = '  +2.22%  '
$ws.Range("D46").Value = '1.710'
$ws.Range("E46").Value = '  +3.56%  '
$ws.Range("D47").Value = '0.06400'
$ws.Range("E47").Value = '  +1.85%  '
$ws.Range("D48").Value = '0.4644'
$ws.Range("E48").Value = '  +1.89%  '
$ws.Range("D49").Value = '1.847'
$ws.Range("E49").Value = '  +2.13%  '
$ws.Range("D50").Value = '39.23'
$ws.Range("E50").Value = '  +4.19%  '
$ws.Range("D51").Value = '63.97'
$ws.Range("E51").Value = '  +0.40%  '

# Restore default (General) styling on these cells so no stray
# cell-level style indices are introduced by the NumberFormat change above.
$ws.Range("D2:E51").ClearFormats()

